# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (GitHub Actions run, Thu May 18 08:32:41 UTC 2023).
# Columns: A=index(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h).
# NumberFormat is forced to Text ("@") before each write so numeric-looking
# price strings (e.g. "313.50", "1.0000") are kept as literal text instead
# of being auto-coerced into numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

    $c = $ws.Cells.Item(2, 4)
    $c.NumberFormat = '@'
    $c.Value = '27.398.23'
    $c = $ws.Cells.Item(2, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +2.15%  '
    $c = $ws.Cells.Item(3, 4)
    $c.NumberFormat = '@'
    $c.Value = '1.827.46'
    $c = $ws.Cells.Item(3, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.04%  '
    $c = $ws.Cells.Item(4, 4)
    $c.NumberFormat = '@'
    $c.Value = '1.0000'
    $c = $ws.Cells.Item(4, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.16%  '
    $c = $ws.Cells.Item(5, 4)
    $c.NumberFormat = '@'
    $c.Value = '313.50'
    $c = $ws.Cells.Item(5, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.18%  '
    $c = $ws.Cells.Item(6, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.12%  '
    $c = $ws.Cells.Item(7, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.4458'
    $c = $ws.Cells.Item(7, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.38%  '
    $c = $ws.Cells.Item(8, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +3.31%  '
    $c = $ws.Cells.Item(9, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.07406'
    $c = $ws.Cells.Item(9, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +2.07%  '
    $c = $ws.Cells.Item(10, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.8816'
    $c = $ws.Cells.Item(10, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +3.66%  '
    $c = $ws.Cells.Item(11, 4)
    $c.NumberFormat = '@'
    $c.Value = '20.91'
    $c = $ws.Cells.Item(11, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.39%  '
    $c = $ws.Cells.Item(12, 4)
    $c.NumberFormat = '@'
    $c.Value = '1.829.08'
    $c = $ws.Cells.Item(12, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.11%  '
    $c = $ws.Cells.Item(13, 4)
    $c.NumberFormat = '@'
    $c.Value = '6.725'
    $c = $ws.Cells.Item(13, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.99%  '
    $c = $ws.Cells.Item(14, 4)
    $c.NumberFormat = '@'
    $c.Value = '5.434'
    $c = $ws.Cells.Item(14, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +2.76%  '
    $c = $ws.Cells.Item(15, 4)
    $c.NumberFormat = '@'
    $c.Value = '92.83'
    $c = $ws.Cells.Item(15, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.95%  '
    $c = $ws.Cells.Item(16, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.20%  '
    $c = $ws.Cells.Item(17, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.13%  '
    $c = $ws.Cells.Item(18, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.000008806'
    $c = $ws.Cells.Item(18, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.25%  '
    $c = $ws.Cells.Item(19, 4)
    $c.NumberFormat = '@'
    $c.Value = '1.000'
    $c = $ws.Cells.Item(19, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.11%  '
    $c = $ws.Cells.Item(20, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.63%  '
    $c = $ws.Cells.Item(21, 4)
    $c.NumberFormat = '@'
    $c.Value = '27.398.44'
    $c = $ws.Cells.Item(21, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +2.05%  '
    $c = $ws.Cells.Item(22, 4)
    $c.NumberFormat = '@'
    $c.Value = '5.369'
    $c = $ws.Cells.Item(22, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +4.62%  '
    $c = $ws.Cells.Item(23, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.30%  '
    $c = $ws.Cells.Item(24, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -1.17%  '
    $c = $ws.Cells.Item(25, 4)
    $c.NumberFormat = '@'
    $c.Value = '151.01'
    $c = $ws.Cells.Item(25, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.13%  '
    $c = $ws.Cells.Item(26, 4)
    $c.NumberFormat = '@'
    $c.Value = '2.298'
    $c = $ws.Cells.Item(26, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +3.67%  '
    $c = $ws.Cells.Item(27, 4)
    $c.NumberFormat = '@'
    $c.Value = '18.63'
    $c = $ws.Cells.Item(27, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.47%  '
    $c = $ws.Cells.Item(28, 4)
    $c.NumberFormat = '@'
    $c.Value = '5.367'
    $c = $ws.Cells.Item(28, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +3.50%  '
    $c = $ws.Cells.Item(29, 4)
    $c.NumberFormat = '@'
    $c.Value = '117.09'
    $c = $ws.Cells.Item(29, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.03%  '
    $c = $ws.Cells.Item(30, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.10%  '
    $c = $ws.Cells.Item(31, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.7954'
    $c = $ws.Cells.Item(31, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +6.74%  '
    $c = $ws.Cells.Item(32, 4)
    $c.NumberFormat = '@'
    $c.Value = '1.200'
    $c = $ws.Cells.Item(32, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +2.62%  '
    $c = $ws.Cells.Item(33, 4)
    $c.NumberFormat = '@'
    $c.Value = '4.575'
    $c = $ws.Cells.Item(33, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +3.60%  '
    $c = $ws.Cells.Item(34, 4)
    $c.NumberFormat = '@'
    $c.Value = '2.932'
    $c = $ws.Cells.Item(34, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +0.08%  '
    $c = $ws.Cells.Item(35, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.9998'
    $c = $ws.Cells.Item(35, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.13%  '
    $c = $ws.Cells.Item(36, 4)
    $c.NumberFormat = '@'
    $c.Value = '1.110'
    $c = $ws.Cells.Item(36, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +2.46%  '
    $c = $ws.Cells.Item(37, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.01985'
    $c = $ws.Cells.Item(37, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.75%  '
    $c = $ws.Cells.Item(38, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.05273'
    $c = $ws.Cells.Item(39, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +3.66%  '
    $c = $ws.Cells.Item(40, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.5327'
    $c = $ws.Cells.Item(40, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.34%  '
    $c = $ws.Cells.Item(41, 4)
    $c.NumberFormat = '@'
    $c.Value = '2.354'
    $c = $ws.Cells.Item(41, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +20.91%  '
    $c = $ws.Cells.Item(42, 4)
    $c.NumberFormat = '@'
    $c.Value = '2.880'
    $c = $ws.Cells.Item(42, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +0.37%  '
    $c = $ws.Cells.Item(43, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.1701'
    $c = $ws.Cells.Item(43, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.04%  '
    $c = $ws.Cells.Item(44, 4)
    $c.NumberFormat = '@'
    $c.Value = '8.676'
    $c = $ws.Cells.Item(44, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +3.15%  '
    $c = $ws.Cells.Item(45, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.5067'
    $c = $ws.Cells.Item(45, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -1.81%  '
    $c = $ws.Cells.Item(46, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.03%  '
    $c = $ws.Cells.Item(47, 2)
    $c.NumberFormat = '@'
    $c.Value = 'PaxosStandard'
    $c = $ws.Cells.Item(47, 3)
    $c.NumberFormat = '@'
    $c.Value = 'https://coinranking.com/coin/B8xT718SbVhhh+paxosstandard-pax'
    $c = $ws.Cells.Item(47, 4)
    $c.NumberFormat = '@'
    $c.Value = '1.001'
    $c = $ws.Cells.Item(47, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.18%  '
    $c = $ws.Cells.Item(48, 2)
    $c.NumberFormat = '@'
    $c.Value = 'Quant'
    $c = $ws.Cells.Item(48, 3)
    $c.NumberFormat = '@'
    $c.Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
    $c = $ws.Cells.Item(48, 4)
    $c.NumberFormat = '@'
    $c.Value = '105.71'
    $c = $ws.Cells.Item(48, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +0.54%  '
    $c = $ws.Cells.Item(49, 2)
    $c.NumberFormat = '@'
    $c.Value = 'NEARProtocol'
    $c = $ws.Cells.Item(49, 3)
    $c.NumberFormat = '@'
    $c.Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    $c = $ws.Cells.Item(49, 4)
    $c.NumberFormat = '@'
    $c.Value = '1.689'
    $c = $ws.Cells.Item(49, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +2.38%  '
    $c = $ws.Cells.Item(50, 2)
    $c.NumberFormat = '@'
    $c.Value = 'PaxDollar'
    $c = $ws.Cells.Item(50, 3)
    $c.NumberFormat = '@'
    $c.Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
    $c = $ws.Cells.Item(50, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.9997'
    $c = $ws.Cells.Item(50, 5)
    $c.NumberFormat = '@'
    $c.Value = '  -0.11%  '
    $c = $ws.Cells.Item(51, 2)
    $c.NumberFormat = '@'
    $c.Value = 'Cronos'
    $c = $ws.Cells.Item(51, 3)
    $c.NumberFormat = '@'
    $c.Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
    $c = $ws.Cells.Item(51, 4)
    $c.NumberFormat = '@'
    $c.Value = '0.06391'
    $c = $ws.Cells.Item(51, 5)
    $c.NumberFormat = '@'
    $c.Value = '  +1.18%  '

